$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.036885187800112
$ws.Range("D2").Value = 1.046176383518953
$ws.Range("E2").Value = 1.035747255347281
$ws.Range("F2").Value = 1.053035732050989
$ws.Range("I2").Value = 1.034564385524929
$ws.Range("J2").Value = 1.041991068914332
$ws.Range("K2").Value = 1.04894255285619
$ws.Range("L2").Value = 1.038542974423473
$ws.Range("M2").Value = 1.055782812123899
$ws.Range("N2").Value = 1.043470815750236

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.037990150260277
$ws.Range("D3").Value = 1.04725260269285
$ws.Range("E3").Value = 1.036689582796459
$ws.Range("F3").Value = 1.054260887623274
$ws.Range("I3").Value = 1.034751415274852
$ws.Range("J3").Value = 1.042739399093332
$ws.Range("K3").Value = 1.049829673808202
$ws.Range("L3").Value = 1.039294407664711
$ws.Range("M3").Value = 1.056819875976903
$ws.Range("N3").Value = 1.044220208643925

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038704639402143
$ws.Range("D4").Value = 1.047948828407883
$ws.Range("E4").Value = 1.037299285452897
$ws.Range("F4").Value = 1.055053781990452
$ws.Range("I4").Value = 1.034870098712282
$ws.Range("J4").Value = 1.043222578670478
$ws.Range("K4").Value = 1.050402946617055
$ws.Range("L4").Value = 1.039779968825536
$ws.Range("M4").Value = 1.057490497739609
$ws.Range("N4").Value = 1.044704074391492

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039004892961898
$ws.Range("D5").Value = 1.048241483878764
$ws.Range("E5").Value = 1.037555593508824
$ws.Range("F5").Value = 1.055387148378964
$ws.Range("I5").Value = 1.034919433423477
$ws.Range("J5").Value = 1.043425458651465
$ws.Range("K5").Value = 1.050643770470399
$ws.Range("L5").Value = 1.039983939398948
$ws.Range("M5").Value = 1.057772325321326
$ws.Range("N5").Value = 1.044907242485323

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039055299982286
$ws.Range("D6").Value = 1.048290619816569
$ws.Range("E6").Value = 1.037598628121025
$ws.Range("F6").Value = 1.055443124057098
$ws.Range("I6").Value = 1.03492768410386
$ws.Range("J6").Value = 1.043459508516417
$ws.Range("K6").Value = 1.050684195309617
$ws.Range("L6").Value = 1.040018177617585
$ws.Range("M6").Value = 1.057819639460538
$ws.Range("N6").Value = 1.044941340704989

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038708651866139
$ws.Range("D7").Value = 1.047952739032618
$ws.Range("E7").Value = 1.037302710296547
$ws.Range("F7").Value = 1.055058236315402
$ws.Range("I7").Value = 1.034870760125099
$ws.Range("J7").Value = 1.043225290539659
$ws.Range("K7").Value = 1.050406165224093
$ws.Range("L7").Value = 1.039782694916193
$ws.Range("M7").Value = 1.057494263932821
$ws.Range("N7").Value = 1.044706790111839

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037258718192485
$ws.Range("D8").Value = 1.046540130271603
$ws.Range("E8").Value = 1.036065728967626
$ws.Range("F8").Value = 1.053449751304537
$ws.Range("I8").Value = 1.034628076916546
$ws.Range("J8").Value = 1.042244186135129
$ws.Range("K8").Value = 1.049242515863286
$ws.Range("L8").Value = 1.038797062908943
$ws.Range("M8").Value = 1.056133381691467
$ws.Range("N8").Value = 1.04372429242652

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.034699912203093
$ws.Range("D9").Value = 1.044049683893755
$ws.Range("E9").Value = 1.033885649230949
$ws.Range("F9").Value = 1.050616403775127
$ws.Range("I9").Value = 1.034182547095278
$ws.Range("J9").Value = 1.040507375099616
$ws.Range("K9").Value = 1.04718622424691
$ws.Range("L9").Value = 1.037055139754492
$ws.Range("M9").Value = 1.053732026759395
$ws.Range("N9").Value = 1.041985014920121

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.032991397544446
$ws.Range("D10").Value = 1.042388509681125
$ws.Range("E10").Value = 1.032432007229662
$ws.Range("F10").Value = 1.048728125020806
$ws.Range("I10").Value = 1.033873509660423
$ws.Range("J10").Value = 1.039344113071726
$ws.Range("K10").Value = 1.045811437078759
$ws.Range("L10").Value = 1.035890402446686
$ws.Range("M10").Value = 1.052128855501452
$ws.Range("N10").Value = 1.040820100926724

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.032250949379333
$ws.Range("D11").Value = 1.041668985608133
$ws.Range("E11").Value = 1.031802499523358
$ws.Range("F11").Value = 1.04791061222313
$ws.Range("I11").Value = 1.033736844315518
$ws.Range("J11").Value = 1.038839123647446
$ws.Range("K11").Value = 1.045215199736646
$ws.Range("L11").Value = 1.0353852339823
$ws.Range("M11").Value = 1.051434115581312
$ws.Range("N11").Value = 1.040314394359539

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.031975814849908
$ws.Range("D12").Value = 1.041401687637058
$ws.Range("E12").Value = 1.031568661174015
$ws.Range("F12").Value = 1.047606968901183
$ws.Range("I12").Value = 1.03368565256465
$ws.Range("J12").Value = 1.038651353736419
$ws.Range("K12").Value = 1.044993587775687
$ws.Range("L12").Value = 1.035197466766405
$ws.Range("M12").Value = 1.051175973592744
$ws.Range("N12").Value = 1.040126357793701

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.032034836648913
$ws.Range("D13").Value = 1.041459025544385
$ws.Range("E13").Value = 1.031618820806389
$ws.Range("F13").Value = 1.047672100662093
$ws.Range("I13").Value = 1.033696652743491
$ws.Range("J13").Value = 1.03869163982529
$ws.Range("K13").Value = 1.045041130771902
$ws.Range("L13").Value = 1.035237749143775
$ws.Range("M13").Value = 1.051231349773729
$ws.Range("N13").Value = 1.040166701093438

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.032228208703602
$ws.Range("D14").Value = 1.041646891387673
$ws.Range("E14").Value = 1.031783170602859
$ws.Range("F14").Value = 1.047885512625744
$ws.Range("I14").Value = 1.033732621523511
$ws.Range("J14").Value = 1.038823606496614
$ws.Range("K14").Value = 1.045196884142034
$ws.Range("L14").Value = 1.035369715641896
$ws.Range("M14").Value = 1.05141277921716
$ws.Range("N14").Value = 1.040298855172572

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.032347338431604
$ws.Range("D15").Value = 1.041762637108548
$ws.Range("E15").Value = 1.031884430443755
$ws.Range("F15").Value = 1.048017005035839
$ws.Range("I15").Value = 1.033754726342017
$ws.Range("J15").Value = 1.038904889744701
$ws.Range("K15").Value = 1.045292829978033
$ws.Range("L15").Value = 1.035451007956481
$ws.Range("M15").Value = 1.05152455263561
$ws.Range("N15").Value = 1.040380253852195

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.033040524831347
$ws.Range("D16").Value = 1.042436257310676
$ws.Range("E16").Value = 1.032473784044295
$ws.Range("F16").Value = 1.04878238315775
$ws.Range("I16").Value = 1.033882519644522
$ws.Range("J16").Value = 1.039377600337372
$ws.Range("K16").Value = 1.045850987398325
$ws.Range("L16").Value = 1.035923911291928
$ws.Range("M16").Value = 1.052174951291082
$ws.Range("N16").Value = 1.040853635748128

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.03347516729152
$ws.Range("D17").Value = 1.042858740656537
$ws.Range("E17").Value = 1.032843450615496
$ws.Range("F17").Value = 1.049262516977541
$ws.Range("I17").Value = 1.033961917985873
$ws.Range("J17").Value = 1.039673773567252
$ws.Range("K17").Value = 1.04620085110643
$ws.Range("L17").Value = 1.036220328603747
$ws.Range("M17").Value = 1.052582779513871
$ws.Range("N17").Value = 1.041150229577966

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.033728624263264
$ws.Range("D18").Value = 1.043105146429992
$ws.Range("E18").Value = 1.033059064041882
$ws.Range("F18").Value = 1.049542583002232
$ws.Range("I18").Value = 1.034007954700217
$ws.Range("J18").Value = 1.039846402048138
$ws.Range("K18").Value = 1.046404829524319
$ws.Range("L18").Value = 1.036393143805857
$ws.Range("M18").Value = 1.052820605226161
$ws.Range("N18").Value = 1.041323103211096

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.033815035962916
$ws.Range("D19").Value = 1.043189160793214
$ws.Range("E19").Value = 1.03313258150396
$ws.Range("F19").Value = 1.049638080370875
$ws.Range("I19").Value = 1.034023605384884
$ws.Range("J19").Value = 1.039905242826201
$ws.Range("K19").Value = 1.046474365455477
$ws.Range("L19").Value = 1.036452055770206
$ws.Range("M19").Value = 1.052901688591168
$ws.Range("N19").Value = 1.041382027549813

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.033428540756411
$ws.Range("D20").Value = 1.04281341443407
$ws.Range("E20").Value = 1.032803789600601
$ws.Range("F20").Value = 1.0492110019215
$ws.Range("I20").Value = 1.033953427738885
$ws.Range("J20").Value = 1.039642009849591
$ws.Range("K20").Value = 1.046163323476859
$ws.Range("L20").Value = 1.036188534109955
$ws.Range("M20").Value = 1.052539028955898
$ws.Range("N20").Value = 1.041118420752182

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.032171268200869
$ws.Range("D21").Value = 1.041591570544571
$ws.Range("E21").Value = 1.03173477399599
$ws.Range("F21").Value = 1.047822667654117
$ws.Range("I21").Value = 1.033722041433267
$ws.Range("J21").Value = 1.038784750971524
$ws.Range("K21").Value = 1.045151022607032
$ws.Range("L21").Value = 1.03533085825659
$ws.Range("M21").Value = 1.051359355113607
$ws.Range("N21").Value = 1.040259944468181

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.031380196882744
$ws.Range("D22").Value = 1.040823147596005
$ws.Range("E22").Value = 1.031062576654471
$ws.Range("F22").Value = 1.046949866389131
$ws.Range("I22").Value = 1.03357408250691
$ws.Range("J22").Value = 1.038244633484564
$ws.Range("K22").Value = 1.044513721883637
$ws.Range("L22").Value = 1.03479087902466
$ws.Range("M22").Value = 1.050617157047196
$ws.Range("N22").Value = 1.039719059952436

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.031799613667886
$ws.Range("D23").Value = 1.041230522528558
$ws.Range("E23").Value = 1.031418927487163
$ws.Range("F23").Value = 1.04741254586899
$ws.Range("I23").Value = 1.033652753130931
$ws.Range("J23").Value = 1.038531066823846
$ws.Range("K23").Value = 1.044851645723211
$ws.Range("L23").Value = 1.035077201075058
$ws.Range("M23").Value = 1.05101065712682
$ws.Range("N23").Value = 1.040005900059914

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.033449609473638
$ws.Range("D24").Value = 1.042833895469382
$ws.Range("E24").Value = 1.032821710726896
$ws.Range("F24").Value = 1.049234279320076
$ws.Range("I24").Value = 1.033957264965628
$ws.Range("J24").Value = 1.039656362890293
$ws.Range("K24").Value = 1.046180280878655
$ws.Range("L24").Value = 1.03620290092006
$ws.Range("M24").Value = 1.052558798115018
$ws.Range("N24").Value = 1.041132794175849

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.035361886182615
$ws.Range("D25").Value = 1.044693675747841
$ws.Range("E25").Value = 1.034449296127219
$ws.Range("F25").Value = 1.051348778480648
$ws.Range("I25").Value = 1.034299846193155
$ws.Range("J25").Value = 1.040957330327623
$ws.Range("K25").Value = 1.047718514569415
$ws.Range("L25").Value = 1.037506076348309
$ws.Range("M25").Value = 1.054353230978479
$ws.Range("N25").Value = 1.042435609136163
